# Generate Report for Archive
# - Update localization status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (columns E & F, the per-locale status columns) and on
#   each locale detail sheet (zh-cn, de-de -> column C, "Status").
# - Narrow the now-shorter status columns to their new autofit width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$newWidth  = 13.4101845877511

# --- Overview sheet: zh-cn (E) and de-de (F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn detail sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de detail sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
